$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF column (I) for rows 27 through 50 to the new value 21.966
$ws.Range("I27:I50").Value = 21.966
